# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 21 (pushing all the
# subsequent rows down by one, through the former last row 44 -> 45), and
# the sheet's dimension grows from A1:T44 to A1:T45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 21, shifting rows 21..44 down to 22..45.
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with this week's record.
$ws.Cells.Item(21, 1).Value = 7
$ws.Cells.Item(21, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, 3).Value = "Ñuble"
$ws.Cells.Item(21, 4).Value = 45272
$ws.Cells.Item(21, 5).Value = 16
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100107
$ws.Cells.Item(21, 8).Value = "Otros"
$ws.Cells.Item(21, 9).Value = 100107002
$ws.Cells.Item(21, 10).Value = "Chirimoya"
$ws.Cells.Item(21, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 50
$ws.Cells.Item(21, 14).Value = 20000
$ws.Cells.Item(21, 15).Value = 20000
$ws.Cells.Item(21, 16).Value = 20000
$ws.Cells.Item(21, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(21, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 19).Value = 2000
$ws.Cells.Item(21, 20).Value = 10
